# This script updates the "Price" (column D) and "Volume(1h)" (column E)
# figures for the crypto-price table, refreshing each row's quoted
# values to the latest scrape while leaving every other cell untouched.
#
# Column D and E cells in the source sheet are stored as plain text
# (not numbers), e.g. "61.796.39" or "  +0.92%  ". Excel's COM layer
# auto-converts any string that *looks* like a plain number (for
# example "0.998" or "79.70") into a real numeric value as soon as it
# is assigned to Range.Value, which would both change the cell type
# and silently drop formatting such as trailing zeros. To keep those
# values as text we prefix them with a leading apostrophe (Excel's
# classic "treat as text" marker) before assignment, then reset the
# cell's Style back to "Normal" so no stray quote-prefix style is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.796.39'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '2.904.89'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'586.54"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = "'145.89"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.35%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.503"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.83%  '
$ws.Range('D9').Value = '2.901.46'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').Value = "'7.11"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = "'0.149"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.37%  '
$ws.Range('D12').Value = "'0.434"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').Value = "'0.0000235"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.40%  '
$ws.Range('D14').Value = "'32.17"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').Value = '3.385.37'
$ws.Range('D17').Value = '61.732.08'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = "'6.58"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('D19').Value = '2.899.91'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = "'433.48"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('D21').Value = "'13.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').Value = "'0.657"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = "'6.91"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('D24').Value = "'79.70"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.55%  '
$ws.Range('D25').Value = "'10.93"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.52%  '
$ws.Range('D26').Value = "'11.86"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.45%  '
$ws.Range('D27').Value = "'2.09"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').Value = "'1.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = "'7.22"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.98%  '
$ws.Range('D30').Value = "'2.56"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').Value = "'0.0000102"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +20.19%  '
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = "'0.108"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.34%  '
$ws.Range('D34').Value = "'0.999"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').Value = "'25.88"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = "'3.07"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.59%  '
$ws.Range('D38').Value = "'5.49"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('D39').Value = "'49.20"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').Value = "'1.98"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').Value = "'8.35"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('D42').Value = "'0.115"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').Value = "'0.273"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.92%  '
$ws.Range('D44').Value = "'39.17"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.31%  '
$ws.Range('D45').Value = "'135.21"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('D46').Value = '2.684.72'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').Value = "'0.0337"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('D48').Value = "'349.42"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').Value = "'0.104"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('D51').Value = "'22.40"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '
